# "Generate Report for Archive" - refresh the localization status report:
#   - flip the Status column from "Ready for handoff" to "In Translation"
#     on every sheet (Overview + each locale tab)
#   - re-fit the Status column width now that the text is shorter

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"

# --- Overview sheet: Status is duplicated in columns E (zh-cn) and F (de-de) ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus

# New status text is narrower, so the Status columns shrink accordingly.
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5

# --- zh-cn sheet: Status is column C ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("C3").Value = $newStatus
$zhcn.Columns.Item(3).ColumnWidth = 12.5

# --- de-de sheet: Status is column C ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = $newStatus
$dede.Range("C3").Value = $newStatus
$dede.Columns.Item(3).ColumnWidth = 12.5
